$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row to append after the last existing data row (row 81 -> row 82)
$newRow = 82
$lastRow = $newRow - 1

# Copy formatting (style) from the row above so the new row matches
# the existing date-format style used throughout column A
$ws.Range("A$lastRow").Copy() | Out-Null
$ws.Range("A$newRow").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Column A: date serial number (Excel 1900 date system) = 45884 -> 2025-08-15
$ws.Cells.Item($newRow, 1).Value = 45884

# Column B: value
$ws.Cells.Item($newRow, 2).Value = 0.06654624964350926
